$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the table (header + 10 data rows)
$data = @(
    @("Cluster Name", "Activecases"),
    @("3641 Calvary The Regent Mount Waverley", 13),
    @("Barwon Heads Hotel Barwon Heads", 36),
    @("Confirmed Omicron Sircuit Bar Fitzroy", 18),
    @("Confirmed Omicron Variant The Peel Hotel Collingwood", 11),
    @("Hotel Traralgon (Ryan's Hotel) Traralgon", 10),
    @("Melbourne Cricket Ground (MCG)", 54),
    @("Melbourne Stars Big Bash Cricket Team East Melbourne", 24),
    @("St Vincents Hospital Melbourne Emergency Department Fitzroy", 15),
    @("Werribee Mercy Hospital Emergency Department", 12),
    @("Western Health Sunshine Hospital Emergency Department St Albans", 12)
)

# Clear the whole previously used range first so stale shared strings are
# dropped and the new strings get appended fresh, in the order we write them.
$oldUsed = $ws.UsedRange
$oldRows = $oldUsed.Rows.Count
$ws.Range("A1:B" + $oldRows).Clear()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
